$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 42636.592673611114
$ws.Range("A5").NumberFormat = "m/d/yy h:mm"
$ws.Range("B5").Value = $true
$ws.Range("C5").Value = 10115.89
$ws.Range("D5").Value = 10085.129999999999
$ws.Range("E5").Value = 81.97
$ws.Range("F5").Value = 81.47
$ws.Range("G5").Value = $true
$ws.Range("G5").NumberFormat = "m/d/yy h:mm"
$ws.Range("H5").Value = -0.61
$ws.Range("I5").Value = $false
